{"js": "// The diff collapses four \"<id>...</id>\" sequences (each currently split\n// across three runs: \"<id>\", the bare id value, \"</id>\") into a single run\n// containing the full \"<id>p112r_N</id>\" text, and fixes a typo\n// (\"mentionned.\" -> \"mentioned.\").\n//\n// Re-inserting the same (now merged) text over each found range keeps the\n// formatting of the first run in the range (matches the target XML, which\n// keeps the <id>/</id> run's rPr: Courier New, color 7f6000, sz/szCs 18).\n\nconst ids = [\"p112r_1\", \"p112r_2\", \"p112r_3\", \"p112r_4\"];\n\nfor (const id of ids) {\n  const searchText = `<id>${id}</id>`;\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find \"${searchText}\" in document body`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(searchText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Fix the \"mentionned.\" -> \"mentioned.\" typo.\nconst typoResults = context.document.body.search(\"mentionned.\", { matchCase: true });\ntypoResults.load(\"text\");\nawait context.sync();\n\nif (typoResults.items.length === 0) {\n  throw new Error('Could not find \"mentionned.\" in document body');\n}\n\nfor (const range of typoResults.items) {\n  range.insertText(\"mentioned.\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The diff collapses four \"<id>...</id>\" sequences (each currently split\n# across three runs: \"<id>\", the bare id value, \"</id>\") into a single run\n# containing the full \"<id>p112r_N</id>\" text, and fixes a typo\n# (\"mentionned.\" -> \"mentioned.\").\n#\n# Find/Replace across the run boundary merges the matched text into a run\n# that keeps the formatting of the first run found (Courier New, color\n# 7f6000, sz/szCs 18) -- matching the target XML.\n\n$d = $word.ActiveDocument\n\n$ids = @(\"p112r_1\", \"p112r_2\", \"p112r_3\", \"p112r_4\")\n\nforeach ($id in $ids) {\n    $searchText = \"<id>$id</id>\"\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $searchText\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"mentionned.\"\n$find.Replacement.Text = \"mentioned.\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n"}
